$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-31 Saturday", "2024-09-01 Sunday"),
    @("612×3=", "120×3="),
    @("528×7=", "386×3="),
    @("653×2=", "382×8="),
    @("469×3=", "838×4="),
    @("143×5=", "476×8="),
    @("503×4=", "732×6="),
    @("562×9=", "530×8="),
    @("376×9=", "888×9="),
    @("774×4=", "209×2="),
    @("982×5=", "563×5="),
    @("483×5=", "856×3="),
    @("110×3=", "161×9="),
    @("966×6=", "415×6="),
    @("113×9=", "225×4="),
    @("660×8=", "843×2="),
    @("217×6=", "355×2="),
    @("697×2=", "440×5="),
    @("848×7=", "643×3="),
    @("538×2=", "434×8="),
    @("848×3=", "497×2="),
    @("513×4=", "937×2="),
    @("700×7=", "240×5="),
    @("945×4=", "340×9="),
    @("499×8=", "901×4="),
    @("905×9=", "832×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
